# Update '想去人数' (F column) values across sheets per commit 456a3b4
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 405
$ws.Range("F7").Value = 115
$ws.Range("F8").Value = 10339
$ws.Range("F10").Value = 3572
$ws.Range("F12").Value = 2462
$ws.Range("F14").Value = 2855
$ws.Range("F17").Value = 2194
$ws.Range("F20").Value = 33
$ws.Range("F21").Value = 399
$ws.Range("F22").Value = 23
$ws.Range("F24").Value = 320
$ws.Range("F26").Value = 242
$ws.Range("F29").Value = 18
$ws.Range("F30").Value = 1263
$ws.Range("F34").Value = 3900
$ws.Range("F35").Value = 3308
$ws.Range("F36").Value = 37
$ws.Range("F38").Value = 1050
$ws.Range("F41").Value = 1296
$ws.Range("F42").Value = 113
$ws.Range("F43").Value = 112
$ws.Range("F44").Value = 74
$ws.Range("F47").Value = 16

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 757
$ws.Range("F3").Value = 997
$ws.Range("F5").Value = 2086

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 757
$ws.Range("F4").Value = 997
$ws.Range("F6").Value = 405
$ws.Range("F11").Value = 115
$ws.Range("F12").Value = 10339
$ws.Range("F15").Value = 3572
$ws.Range("F16").Value = 2462
$ws.Range("F18").Value = 2855
$ws.Range("F20").Value = 2194
$ws.Range("F23").Value = 33
$ws.Range("F24").Value = 399
$ws.Range("F25").Value = 23
$ws.Range("F26").Value = 320
$ws.Range("F28").Value = 242
$ws.Range("F31").Value = 18
$ws.Range("F32").Value = 1263
$ws.Range("F36").Value = 3900
$ws.Range("F37").Value = 3308
$ws.Range("F38").Value = 37
$ws.Range("F39").Value = 1050
$ws.Range("F44").Value = 1296
$ws.Range("F45").Value = 113
$ws.Range("F46").Value = 74
$ws.Range("F48").Value = 16
